$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.396.62"
$ws.Range("E2").Value = "'  +1.85%  "
$ws.Range("D3").Value = "'3.608.06"
$ws.Range("E3").Value = "'  +0.68%  "
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("D5").Value = "'207.04"
$ws.Range("E5").Value = "'  +7.70%  "
$ws.Range("D6").Value = "'571.90"
$ws.Range("E6").Value = "'  -0.50%  "
$ws.Range("E7").Value = "'  -0.31%  "
$ws.Range("E9").Value = "'  +0.97%  "
$ws.Range("D10").Value = "'64.35"
$ws.Range("E10").Value = "'  +15.14%  "
$ws.Range("D11").Value = "'0.149"
$ws.Range("E11").Value = "'  -0.25%  "
$ws.Range("E12").Value = "'  +4.30%  "
$ws.Range("D13").Value = "'10.31"
$ws.Range("E13").Value = "'  +4.56%  "
$ws.Range("D14").Value = "'4.183.97"
$ws.Range("E14").Value = "'  +0.63%  "
$ws.Range("D15").Value = "'3.627.88"
$ws.Range("E15").Value = "'  +1.17%  "
$ws.Range("D16").Value = "'19.33"
$ws.Range("E16").Value = "'  +5.05%  "
$ws.Range("E17").Value = "'  +0.41%  "
$ws.Range("D18").Value = "'68.233.88"
$ws.Range("E18").Value = "'  +1.71%  "
$ws.Range("D19").Value = "'12.31"
$ws.Range("E19").Value = "'  +1.05%  "
$ws.Range("E20").Value = "'  +0.97%  "
$ws.Range("D21").Value = "'406.12"
$ws.Range("E21").Value = "'  +1.34%  "
$ws.Range("D22").Value = "'4.20"
$ws.Range("E22").Value = "'  +0.06%  "
$ws.Range("D23").Value = "'12.47"
$ws.Range("E23").Value = "'  +9.54%  "
$ws.Range("D24").Value = "'85.12"
$ws.Range("E24").Value = "'  -1.03%  "
$ws.Range("D25").Value = "'2.92"
$ws.Range("E25").Value = "'  -0.67%  "
$ws.Range("E26").Value = "'  +0.87%  "
$ws.Range("D27").Value = "'3.86"
$ws.Range("E27").Value = "'  +6.90%  "
$ws.Range("E28").Value = "'  +3.41%  "
$ws.Range("D29").Value = "'7.77"
$ws.Range("E29").Value = "'  +1.37%  "
$ws.Range("D30").Value = "'31.74"
$ws.Range("E30").Value = "'  +1.60%  "
$ws.Range("D31").Value = "'717.39"
$ws.Range("E31").Value = "'  +12.48%  "
$ws.Range("E32").Value = "'  +0.55%  "
$ws.Range("E33").Value = "'  +0.11%  "
$ws.Range("D34").Value = "'63.71"
$ws.Range("E34").Value = "'  -0.55%  "
$ws.Range("D35").Value = "'42.28"
$ws.Range("E35").Value = "'  -0.09%  "
$ws.Range("D36").Value = "'0.421"
$ws.Range("E36").Value = "'  +5.26%  "
$ws.Range("E37").Value = "'  -0.04%  "
$ws.Range("D38").Value = "'3.31"
$ws.Range("E38").Value = "'  +11.58%  "
$ws.Range("D39").Value = "'0.0₃0756"
$ws.Range("E39").Value = "'  -2.27%  "
$ws.Range("D40").Value = "'3.14"
$ws.Range("E40").Value = "'  +20.11%  "
$ws.Range("D41").Value = "'3.196.43"
$ws.Range("E41").Value = "'  -0.13%  "
$ws.Range("E42").Value = "'  -0.47%  "
$ws.Range("E43").Value = "'  +0.02%  "
$ws.Range("D44").Value = "'2.68"
$ws.Range("E44").Value = "'  -0.37%  "
$ws.Range("E45").Value = "'  +9.90%  "
$ws.Range("D46").Value = "'0.0417"
$ws.Range("E46").Value = "'  +0.10%  "
$ws.Range("D47").Value = "'0.132"
$ws.Range("E47").Value = "'  +0.87%  "
$ws.Range("D48").Value = "'8.82"
$ws.Range("E48").Value = "'  +2.60%  "
$ws.Range("E49").Value = "'  -1.02%  "
$ws.Range("E50").Value = "'  -2.33%  "
$ws.Range("E51").Value = "'  -0.94%  "
